$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 33 - this shifts the existing rows 33-90
# down to 34-91 (matching the new D/J/K/L/M/N/O/P/Q values seen throughout
# the diff, which is just every record sliding down by one row).
$ws.Rows("33:33").Insert()

# Populate the newly inserted row 33 with the new weekly price record.
$ws.Range("A33").Value = 9
$ws.Range("B33").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C33").Value = "Metropolitana"
$ws.Range("D33").Value = 45219
$ws.Range("E33").Value = 13
$ws.Range("F33").Value = 100112010
$ws.Range("G33").Value = "Achicoria"
$ws.Range("H33").Value = "Sin especificar"
$ws.Range("I33").Value = "Primera"
$ws.Range("J33").Value = 70
$ws.Range("K33").Value = 7000
$ws.Range("L33").Value = 7000
$ws.Range("M33").Value = 7000
$ws.Range("N33").Value = "$/caja 16 unidades"
$ws.Range("O33").Value = "Provincia de Quillota"
$ws.Range("P33").Value = 438
$ws.Range("Q33").Value = 16
$ws.Range("R33").Value = "Hortaliza"
